# Auto-generated edit script: updates crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '60.912.82'
$ws.Range("E2").Value = '  -1.37%  '

# Row 3
$ws.Range("D3").Value = '3.409.97'
$ws.Range("E3").Value = '  -1.15%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.29%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.48'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.40%  '

# Row 7
$ws.Range("E7").Value = '  +0.06%  '

# Row 8
$ws.Range("D8").Value = '3.411.06'
$ws.Range("E8").Value = '  -1.11%  '

# Row 9
$ws.Range("E9").Value = '  +1.20%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.54'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.77%  '

# Row 11
$ws.Range("E11").Value = '  +1.20%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.392'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.07%  '

# Row 13
$ws.Range("D13").Value = '3.991.81'
$ws.Range("E13").Value = '  -0.98%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.21'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.95%  '

# Row 15
$ws.Range("E15").Value = '  +0.53%  '

# Row 16
$ws.Range("E16").Value = '  -1.11%  '

# Row 17
$ws.Range("D17").Value = '3.420.02'
$ws.Range("E17").Value = '  -0.97%  '

# Row 18
$ws.Range("D18").Value = '61.011.38'
$ws.Range("E18").Value = '  -1.24%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.31'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.10%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.37'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.59%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.31%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '395.34'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.49%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.566'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.23%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.05'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.01%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.996'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.68%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000123'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.66%  '

# Row 27
$ws.Range("D27").Value = '3.557.82'
$ws.Range("E27").Value = '  -0.91%  '

# Row 28
$ws.Range("E28").Value = '  +0.70%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.45'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.91%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.12%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.85%  '

# Row 32
$ws.Range("E32").Value = '  -0.24%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.44'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.16%  '

# Row 34
$ws.Range("E34").Value = '  -0.01%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.87'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.33%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.38%  '

# Row 37
$ws.Range("D37").Value = '3.435.88'
$ws.Range("E37").Value = '  -1.04%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.14'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.42%  '

# Row 39
$ws.Range("E39").Value = '  -0.87%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '167.59'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.80%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0786'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.34%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.01'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.54%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.795'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.15%  '

# Row 44
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.51'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.89%  '

# Row 45
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.10%  '

# Row 46
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.92'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.43%  '

# Row 47
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.71'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.83%  '

# Row 48
$ws.Range("D48").Value = '2.584.11'
$ws.Range("E48").Value = '  -1.15%  '

# Row 49
$ws.Range("E49").Value = '  -4.54%  '

# Row 50
$ws.Range("E50").Value = '  +1.23%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.96'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.13%  '
